$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Price" (D) and "Volume(1h)" (E) columns for rows 2-51
# with refreshed crypto data. All of these cells are plain text in the
# workbook (e.g. "25.904.02", "  -2.11%  "), so for any new value that
# would otherwise be auto-parsed by Excel as a number (like "215.17"),
# a leading apostrophe forces it to stay text; the cell style is then
# reset to "Normal" so no stray number-format/style gets attached.

$ws.Range("D2").Value = "25.904.02"
$ws.Range("E2").Value = "  -2.11%  "
$ws.Range("D3").Value = "1.632.86"
$ws.Range("E3").Value = "  -2.21%  "
$ws.Range("D4").Value = "'1.016"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.85%  "
$ws.Range("D5").Value = "'215.17"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.68%  "
$ws.Range("D6").Value = "'0.5030"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.32%  "
$ws.Range("D7").Value = "'1.016"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.93%  "
$ws.Range("D8").Value = "'0.2575"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("D9").Value = "'0.06393"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.17%  "
$ws.Range("D10").Value = "'19.45"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.10%  "
$ws.Range("D11").Value = "'0.07763"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.23%  "
$ws.Range("D12").Value = "1.638.97"
$ws.Range("E12").Value = "  -2.07%  "
$ws.Range("D13").Value = "'4.259"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.00%  "
$ws.Range("D14").Value = "1.859.87"
$ws.Range("E14").Value = "  -2.16%  "
$ws.Range("D15").Value = "'0.5430"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.58%  "
$ws.Range("D16").Value = "0.0₅7933"
$ws.Range("E16").Value = "  -1.60%  "
$ws.Range("D17").Value = "'63.43"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.12%  "
$ws.Range("D18").Value = "25.915.49"
$ws.Range("E18").Value = "  -2.28%  "
$ws.Range("D19").Value = "'1.017"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.01%  "
$ws.Range("D20").Value = "'203.60"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.32%  "
$ws.Range("D21").Value = "'4.302"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.04%  "
$ws.Range("D22").Value = "'9.981"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.47%  "
$ws.Range("D23").Value = "'5.962"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.97%  "
$ws.Range("D24").Value = "'1.018"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.05%  "
$ws.Range("D25").Value = "'1.974"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +13.51%  "
$ws.Range("D26").Value = "'141.80"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.37%  "
$ws.Range("D27").Value = "'0.1151"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.29%  "
$ws.Range("D28").Value = "'15.68"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.59%  "
$ws.Range("D29").Value = "'6.793"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.16%  "
$ws.Range("D30").Value = "'1.240"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.71%  "
$ws.Range("D31").Value = "'0.04984"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.55%  "
$ws.Range("D32").Value = "'3.255"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.55%  "
$ws.Range("D33").Value = "'3.189"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.10%  "
$ws.Range("D34").Value = "'1.535"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.06%  "
$ws.Range("D35").Value = "'2.350"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.12%  "
$ws.Range("D36").Value = "'2.624"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.27%  "
$ws.Range("D37").Value = "'0.8889"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.92%  "
$ws.Range("D38").Value = "'0.5640"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.92%  "
$ws.Range("D39").Value = "1.111.57"
$ws.Range("E39").Value = "  -4.75%  "
$ws.Range("D40").Value = "'0.01564"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.37%  "
$ws.Range("D41").Value = "'2.588"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.42%  "
$ws.Range("D42").Value = "'1.017"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.02%  "
$ws.Range("D43").Value = "'5.611"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.59%  "
$ws.Range("D44").Value = "'0.8154"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.49%  "
$ws.Range("D45").Value = "'99.57"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.84%  "
$ws.Range("D46").Value = "1.772.39"
$ws.Range("E46").Value = "  -2.10%  "
$ws.Range("D47").Value = "0.0₈112"
$ws.Range("E47").Value = "  +1.33%  "
$ws.Range("D48").Value = "'0.4554"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.35%  "
$ws.Range("D49").Value = "'1.019"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.23%  "
$ws.Range("D50").Value = "'54.69"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.61%  "
$ws.Range("D51").Value = "'0.05034"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.76%  "
